{"js": "// Remove the \"Due: ...\" paragraph and the \"This assignment is pass/fail...\"\n// paragraph, per the commit \"Latest slides and programs\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [];\nfor (const p of paragraphs.items) {\n  const t = p.text.trim();\n  if (t.startsWith(\"Due:\") && t.includes(\"23:59pm\")) {\n    targets.push(p);\n  } else if (t.startsWith(\"This assignment is pass/fail\")) {\n    targets.push(p);\n  }\n}\n\nfor (const p of targets) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Due: ...\" paragraph and the \"This assignment is pass/fail...\"\n# paragraph, per the commit \"Latest slides and programs\".\n$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $para = $d.Paragraphs.Item($i)\n    $t = $para.Range.Text.Trim()\n    if ($t.StartsWith(\"Due:\") -and $t.Contains(\"23:59pm\")) {\n        $para.Range.Delete()\n    }\n    elseif ($t.StartsWith(\"This assignment is pass/fail\")) {\n        $para.Range.Delete()\n    }\n}\n"}
